$d = $word.ActiveDocument

$p = $d.Paragraphs(1)
$insertPoint = $p.Range
$insertPoint.Collapse(0)
$insertPoint.InsertAfter(" - 6666")

$newTextRange = $d.Range(5, 12)
$newTextRange.LanguageID = "en-US"

$paraMarkRange = $d.Range(12, 13)
$paraMarkRange.LanguageID = "en-US"
